$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Neo4j query in B2 to add a `Description` field derived from
# sa.arm_descriptions instead of being produced directly in the RETURN clause.
$newQuery = "MATCH (s:study)<-[:member_of]-(sa:study_arm)`nWHERE s.clinical_study_designation in ['COTC007B']`nOPTIONAL MATCH (sa)<-[:member_of]-(c:cohort)`nWITH c, sa`nORDER BY toInteger(left(c.cohort_dose, size(c.cohort_dose) - 9))`nWITH`n    sa,`n    coalesce(sa.arm) as ``Arms``,`n   coalesce(sa.arm_descriptions,`"`") as ``Description``,`n    COLLECT(DISTINCT c.cohort_dose) as cohortDoses`nRETURN`n``Arms``,`n``Description``,`nREDUCE(s = `"`", dose IN cohortDoses | s + dose ) as Cohorts`norder by ``Arms``  asc"
$ws.Range("B2").Value = $newQuery

# Remove the obsolete StatQuery column (column C) entirely.
$ws.Columns.Item(3).Delete()

# Remove the stray leftover formatted row far below the data.
$ws.Rows.Item(21).Delete()

# The header row no longer has the oversized StatQuery font, so it shrinks
# back down to the sheet's default row height.
$ws.Rows.Item(1).AutoFit()

# Column A widened/auto-fit now that the sheet only holds the remaining columns.
$ws.Columns.Item(1).ColumnWidth = 21.6

$ws.Columns.Item(3).Select()
